# Auto-generated script to apply numeric cell updates across multiple sheets
# in the Asura_Profits workbook, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1696.4366
$ws.Range("I15").Value = 1696.4366
$ws.Range("K15").Value = 5089.3098
$ws.Range("M15").Value = -4920.3098
$ws.Range("H121").Value = 1549.1666
$ws.Range("J121").Value = 1713.3334
$ws.Range("L121").Value = 5140.0002
$ws.Range("N121").Value = -8634.0002
$ws.Range("H130").Value = 45833.332
$ws.Range("J130").Value = 45833.332
$ws.Range("L130").Value = 45833.332
$ws.Range("N130").Value = -55873.332
$ws.Range("H132").Value = 2063.1475
$ws.Range("I132").Value = 1602.3043
$ws.Range("J132").Value = 3476.4
$ws.Range("K132").Value = 4806.9129
$ws.Range("L132").Value = 10429.2
$ws.Range("M132").Value = -2276.9129
$ws.Range("N132").Value = -15489.2
$ws.Range("H137").Value = 1487.4546
$ws.Range("I137").Value = 1241.875
$ws.Range("K137").Value = 3725.625
$ws.Range("M137").Value = -1175.625
$ws.Range("H138").Value = 2958.078
$ws.Range("I138").Value = 1967.186
$ws.Range("J138").Value = 4211.2646
$ws.Range("K138").Value = 5901.558
$ws.Range("L138").Value = 12633.7938
$ws.Range("M138").Value = -761.558
$ws.Range("N138").Value = -22913.7938

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 24853.818
$ws.Range("J123").Value = 24853.818
$ws.Range("L123").Value = 24853.818
$ws.Range("N123").Value = -34653.818

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3365666.8
$ws.Range("I7").Value = 10000000
$ws.Range("K7").Value = 10000000
$ws.Range("M7").Value = -9999887
$ws.Range("H35").Value = 15250
$ws.Range("J35").Value = 15250
$ws.Range("L35").Value = 15250
$ws.Range("N35").Value = -15870
$ws.Range("H134").Value = 2112.1956
$ws.Range("I134").Value = 1621.25
$ws.Range("J134").Value = 3234.3572
$ws.Range("K134").Value = 4863.75
$ws.Range("L134").Value = 9703.071599999999
$ws.Range("M134").Value = -2328.75
$ws.Range("N134").Value = -14773.0716

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1677.4375
$ws.Range("I31").Value = 1231.1666
$ws.Range("J31").Value = 4801.3335
$ws.Range("K31").Value = 1231.1666
$ws.Range("L31").Value = 4801.3335
$ws.Range("M31").Value = -936.1666
$ws.Range("N31").Value = -5391.3335
$ws.Range("H34").Value = 1677.4375
$ws.Range("I34").Value = 1231.1666
$ws.Range("J34").Value = 4801.3335
$ws.Range("K34").Value = 1231.1666
$ws.Range("L34").Value = 4801.3335
$ws.Range("M34").Value = -1029.1666
$ws.Range("N34").Value = -5205.3335
$ws.Range("H58").Value = 904691.25
$ws.Range("I58").Value = 1611360.2
$ws.Range("J58").Value = 1725.3334
$ws.Range("K58").Value = 1611360.2
$ws.Range("L58").Value = 1725.3334
$ws.Range("M58").Value = -1611157.2
$ws.Range("N58").Value = -2131.3334
$ws.Range("H132").Value = 752825.1
$ws.Range("I132").Value = 966917
$ws.Range("J132").Value = 3503.5
$ws.Range("K132").Value = 2900751
$ws.Range("L132").Value = 10510.5
$ws.Range("M132").Value = -2898221
$ws.Range("N132").Value = -15570.5
$ws.Range("H136").Value = 904691.25
$ws.Range("I136").Value = 1611360.2
$ws.Range("J136").Value = 1725.3334
$ws.Range("K136").Value = 4834080.6
$ws.Range("L136").Value = 5176.0002
$ws.Range("M136").Value = -4831530.6
$ws.Range("N136").Value = -10276.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6566
$ws.Range("H33").Value = 1143.4166
$ws.Range("J33").Value = 2076.6
$ws.Range("L33").Value = 12459.6
$ws.Range("N33").Value = -13025.6
$ws.Range("H98").Value = 1749.8889
$ws.Range("I98").Value = 5150
$ws.Range("J98").Value = 778.4286
$ws.Range("K98").Value = 15450
$ws.Range("L98").Value = 2335.2858
$ws.Range("M98").Value = -13952
$ws.Range("N98").Value = -5331.2858
$ws.Range("H124").Value = 13000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 13000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39000
$ws.Range("N124").Value = -48820
$ws.Range("H125").Value = 3852.8572
$ws.Range("I125").Value = 2015
$ws.Range("J125").Value = 4588
$ws.Range("K125").Value = 6045
$ws.Range("L125").Value = 13764
$ws.Range("M125").Value = -1125
$ws.Range("N125").Value = -23604
$ws.Range("H126").Value = 3790
$ws.Range("I126").Value = 3030
$ws.Range("J126").Value = 3837.5
$ws.Range("K126").Value = 9090
$ws.Range("L126").Value = 11512.5
$ws.Range("M126").Value = -4150
$ws.Range("N126").Value = -21392.5
$ws.Range("H129").Value = 4167742.8
$ws.Range("I129").Value = 993.3333
$ws.Range("J129").Value = 5556659
$ws.Range("K129").Value = 2979.9999
$ws.Range("L129").Value = 16669977
$ws.Range("M129").Value = 2020.0001
$ws.Range("N129").Value = -16679977
$ws.Range("H131").Value = 2423.707
$ws.Range("I131").Value = 9429.166999999999
$ws.Range("J131").Value = 1457.4368
$ws.Range("K131").Value = 28287.501
$ws.Range("L131").Value = 4372.3104
$ws.Range("M131").Value = -23247.501
$ws.Range("N131").Value = -14452.3104
$ws.Range("M124").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 18819.53
$ws.Range("J123").Value = 18819.53
$ws.Range("L123").Value = 18819.53
$ws.Range("N123").Value = -23719.53

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6361.6665
$ws.Range("I9").Value = 438
$ws.Range("J9").Value = 24132.666
$ws.Range("K9").Value = 438
$ws.Range("L9").Value = 24132.666
$ws.Range("M9").Value = -214
$ws.Range("N9").Value = -24580.666
$ws.Range("H18").Value = 12663
$ws.Range("J18").Value = 12663
$ws.Range("L18").Value = 12663
$ws.Range("N18").Value = -13007
$ws.Range("H20").Value = 30025000
$ws.Range("I20").Value = 30025000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 30025000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -30024774
$ws.Range("H132").Value = 2923.0527
$ws.Range("I132").Value = 2675.7334
$ws.Range("J132").Value = 3850.5
$ws.Range("K132").Value = 8027.2002
$ws.Range("L132").Value = 11551.5
$ws.Range("M132").Value = -5497.2002
$ws.Range("N132").Value = -16611.5
$ws.Range("N20").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 22863.143
$ws.Range("J15").Value = 22863.143
$ws.Range("L15").Value = 22863.143
$ws.Range("N15").Value = -23439.143
$ws.Range("H54").Value = 45500
$ws.Range("J54").Value = 45500
$ws.Range("L54").Value = 45500
$ws.Range("N54").Value = -46540
$ws.Range("H81").Value = 334330.34
$ws.Range("I81").Value = 334330.34
$ws.Range("K81").Value = 668660.6800000001
$ws.Range("M81").Value = -667599.6800000001
$ws.Range("H84").Value = 334330.34
$ws.Range("I84").Value = 334330.34
$ws.Range("K84").Value = 3343303.4
$ws.Range("M84").Value = -3337999.4
$ws.Range("H123").Value = 22207.37
$ws.Range("J123").Value = 22207.37
$ws.Range("L123").Value = 22207.37
$ws.Range("N123").Value = -32007.37
$ws.Range("H132").Value = 1233.4482
$ws.Range("I132").Value = 1050.2273
$ws.Range("J132").Value = 1809.2858
$ws.Range("K132").Value = 3150.6819
$ws.Range("L132").Value = 5427.857400000001
$ws.Range("M132").Value = -620.6819
$ws.Range("N132").Value = -10487.8574

Write-Host "Applied all cell updates."